$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fulano -> editor, with new password hash + new hash value
$ws.Range("B3").Value = 'editor'
$ws.Range("C3").Value = '$2b$12$v1uTTh7BSg4Y9IxSekCoXORKtLD1CCeoyLC6MeWuDWjBkCQQUuDfuP/54UjRv6cR6.Og.haFd8dFi6q0z5t77dCy7iQCndk2XG'
$ws.Range("D3").Value = '744f26fa641bc48221956243d43d6a3dfcf88776927f3a5cc397eb130d910630'

# Row 4: beltrano -> team, with new password hash + new hash value
$ws.Range("B4").Value = 'team'
$ws.Range("C4").Value = '$2b$12$e5xG2DT48B0fecxMI7Qhe..f1D47.vXuxvxCLlmHkePwuurFvf3nC'
$ws.Range("D4").Value = '744f26fa641bc48221956243d43d6a3dfcf88776927f3a5cc397eb130d910630'

# Style updates: E2, E3, E4 become centered (matches header style); D4 becomes underlined to match D3
$ws.Range("E2:E4").HorizontalAlignment = -4108
$ws.Range("E2:E4").Font.Size = 12

$ws.Range("D4").Font.Underline = $true
$ws.Range("D4").Font.Size = 12
$ws.Range("D4").HorizontalAlignment = -4131

# Column E default style switches from left-aligned to centered
$ws.Columns("E").HorizontalAlignment = -4108
$ws.Columns("E").Font.Size = 12

# Update the active selection cell to D5
$ws.Range("D5").Select()
